$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (column G) for rows 2 and 3
$wsOverview.Range("G2").Value = "2017-02-17 09:58:33"
$wsOverview.Range("G3").Value = "2017-02-17 09:58:33"

# zh-cn sheet: Priority (E) ht -> mt, Correspond Handoff Datetime (H), Correspond Handback DateTime (L)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2017-02-17 09:58:16"
$wsZhCn.Range("H3").Value = "2017-02-17 09:58:16"
$wsZhCn.Range("L2").Value = "2017-02-17 09:59:08"
$wsZhCn.Range("L3").Value = "2017-02-17 09:59:08"

# de-de sheet: Correspond Handback DateTime (L)
$wsDeDe.Range("L2").Value = "2017-02-17 09:59:33"
$wsDeDe.Range("L3").Value = "2017-02-17 09:59:33"
